$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Correct a handful of DBH (diameter at breast height) measurement errors ---
$ws.Range("D21").Value = 33.6
$ws.Range("D23").Value = 32.5
$ws.Range("D35").Value = 34.4

# --- Flag rows whose DBH was derived from a circumference measurement that is
#     shared with (mapped from) another tree number, by noting it in column K ---
$note = "Note: alternative measurement shared with tree number mapping"
$ws.Range("K22").Value = $note
$ws.Range("K23").Value = $note
$ws.Range("K24").Value = $note
$ws.Range("K28").Value = $note
$ws.Range("K31").Value = $note
$ws.Range("K32").Value = $note
$ws.Range("K33").Value = $note
$ws.Range("K35").Value = $note
$ws.Range("K36").Value = $note
$ws.Range("K41").Value = $note

# --- Re-position the saved view/selection (re-ran script, scrolled up, selected D36) ---
$ws.Activate() | Out-Null
$ws.Range("D36").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 2
